$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue "D2" "69.655.44"
Set-TextValue "E2" "  +0.73%  "
Set-TextValue "D3" "3.484.49"
Set-TextValue "E3" "  -0.49%  "
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "598.86"
Set-TextValue "E5" "  -0.93%  "
Set-TextValue "D6" "172.44"
Set-TextValue "E6" "  +1.66%  "
Set-TextValue "D7" "0.604"
Set-TextValue "E7" "  -0.63%  "
Set-TextValue "D8" "3.483.17"
Set-TextValue "E8" "  -0.33%  "
Set-TextValue "E9" "  +0.00%  "
Set-TextValue "E10" "  -1.87%  "
Set-TextValue "D11" "7.21"
Set-TextValue "E11" "  +6.40%  "
Set-TextValue "D12" "0.577"
Set-TextValue "E12" "  -0.05%  "
Set-TextValue "D13" "45.71"
Set-TextValue "E13" "  -2.95%  "
Set-TextValue "E14" "  -2.04%  "
Set-TextValue "D15" "4.042.60"
Set-TextValue "E15" "  -0.56%  "
Set-TextValue "D16" "609.93"
Set-TextValue "E16" "  -1.05%  "
Set-TextValue "D17" "8.21"
Set-TextValue "E17" "  -1.65%  "
Set-TextValue "D18" "3.497.15"
Set-TextValue "E18" "  -0.27%  "
Set-TextValue "D19" "69.699.87"
Set-TextValue "E19" "  +0.75%  "
Set-TextValue "E20" "  +0.64%  "
Set-TextValue "D21" "17.08"
Set-TextValue "E21" "  -0.73%  "
Set-TextValue "D22" "0.864"
Set-TextValue "E22" "  -1.22%  "
Set-TextValue "D23" "8.85"
Set-TextValue "E23" "  -20.75%  "
Set-TextValue "D24" "15.39"
Set-TextValue "E24" "  -2.49%  "
Set-TextValue "D25" "95.37"
Set-TextValue "E25" "  -0.75%  "
Set-TextValue "D26" "3.67"
Set-TextValue "E26" "  -4.43%  "
Set-TextValue "E27" "  -0.15%  "
Set-TextValue "D28" "2.53"
Set-TextValue "E28" "  -2.93%  "
Set-TextValue "D29" "33.75"
Set-TextValue "E29" "  +1.35%  "
Set-TextValue "D30" "8.82"
Set-TextValue "E30" "  -3.93%  "
Set-TextValue "D31" "8.01"
Set-TextValue "E31" "  -4.92%  "
Set-TextValue "D32" "2.94"
Set-TextValue "E32" "  -4.81%  "
Set-TextValue "E33" "  -4.05%  "
Set-TextValue "B34" "Bittensor"
Set-TextValue "C34" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D34" "629.47"
Set-TextValue "E34" "  +10.50%  "
Set-TextValue "B35" "NEARProtocol"
Set-TextValue "C35" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D35" "6.76"
Set-TextValue "E35" "  -1.55%  "
Set-TextValue "D36" "0.0985"
Set-TextValue "E36" "  -2.76%  "
Set-TextValue "D37" "3.50"
Set-TextValue "E37" "  -0.62%  "
Set-TextValue "D38" "10.61"
Set-TextValue "E38" "  -1.06%  "
Set-TextValue "D39" "0.0469"
Set-TextValue "E39" "  +6.81%  "
Set-TextValue "E40" "  +0.24%  "
Set-TextValue "D41" "56.06"
Set-TextValue "E41" "  -1.85%  "
Set-TextValue "E42" "  +1.79%  "
Set-TextValue "D43" "3.317.53"
Set-TextValue "E43" "  -2.24%  "
Set-TextValue "D44" "0.0₃0723"
Set-TextValue "E44" "  +2.49%  "
Set-TextValue "D45" "0.307"
Set-TextValue "E45" "  -5.51%  "
Set-TextValue "D46" "2.88"
Set-TextValue "E46" "  +2.20%  "
Set-TextValue "D47" "31.68"
Set-TextValue "E47" "  -3.18%  "
Set-TextValue "D48" "2.52"
Set-TextValue "E48" "  -1.79%  "
Set-TextValue "D49" "0.128"
Set-TextValue "E49" "  -0.19%  "
Set-TextValue "D50" "133.69"
Set-TextValue "E50" "  -0.19%  "
Set-TextValue "E51" "  -0.02%  "

Write-Host "Applied 94 cell updates"
